{"js": "// Load the body's paragraphs so we can find the title and the last paragraph.\nconst body = context.document.body;\nconst paragraphs = body.paragraphs;\nparagraphs.load(\"items\");\nawait context.sync();\n\n// --- 1) Title paragraph: drop the \"Balk1\" (Heading 1) style, replace it with\n//        direct bold formatting on both the paragraph mark and the run, by\n//        replacing the paragraph's OOXML outright (keeps everything else,\n//        e.g. w:rsid* attributes, from being re-added). ---\nconst titlePara = paragraphs.items[0];\nconst titleRange = titlePara.getRange(\"Whole\");\nconst titleOoxml =\n  '<?xml version=\"1.0\" encoding=\"UTF-8\" standalone=\"yes\"?>' +\n  '<pkg:package xmlns:pkg=\"http://schemas.microsoft.com/office/2006/xmlPackage\">' +\n  '<pkg:part pkg:name=\"/word/document.xml\" pkg:contentType=\"application/vnd.openxmlformats-officedocument.wordprocessingml.document.main+xml\">' +\n  '<pkg:xmlData>' +\n  '<w:document xmlns:w=\"http://schemas.openxmlformats.org/wordprocessingml/2006/main\">' +\n  '<w:body>' +\n  '<w:p>' +\n  '<w:pPr><w:rPr><w:b/></w:rPr></w:pPr>' +\n  '<w:r><w:rPr><w:b/></w:rPr><w:t>Git Dersi \u00c7al\u0131\u015fmalar\u0131</w:t></w:r>' +\n  '</w:p>' +\n  '</w:body>' +\n  '</w:document>' +\n  '</pkg:xmlData>' +\n  '</pkg:part>' +\n  '</pkg:package>';\ntitleRange.insertOoxml(titleOoxml, \"Replace\");\nawait context.sync();\n\n// --- 2) Append two new paragraphs after the existing \"Belgemin ilk\n//        paragraf\u0131\" paragraph: \"\u0130kinci paragraf\u0131\" and \"\u00dc\u00e7\u00fcnc\u00fc paragraf\u0131\". ---\nparagraphs.load(\"items\");\nawait context.sync();\nconst lastPara = paragraphs.items[paragraphs.items.length - 1];\nconst secondPara = lastPara.insertParagraph(\"\u0130kinci paragraf\u0131\", \"After\");\nawait context.sync();\nsecondPara.insertParagraph(\"\u00dc\u00e7\u00fcnc\u00fc paragraf\u0131\", \"After\");\nawait context.sync();\n", "ps1": "$d = $word.ActiveDocument\n\n# --- 1) Title paragraph: drop the \"Balk1\" (Heading 1) style, replace it with\n#        direct bold formatting on both the paragraph mark and the run, by\n#        replacing that paragraph's range contents via OOXML (InsertXML\n#        replaces only the targeted range, leaving the rest of the document\n#        untouched). ---\n$titlePara = $d.Paragraphs(1)\n$titleRange = $titlePara.Range\n$titleOoxml = '<?xml version=\"1.0\" encoding=\"UTF-8\" standalone=\"yes\"?><pkg:package xmlns:pkg=\"http://schemas.microsoft.com/office/2006/xmlPackage\"><pkg:part pkg:name=\"/word/document.xml\" pkg:contentType=\"application/vnd.openxmlformats-officedocument.wordprocessingml.document.main+xml\"><pkg:xmlData><w:document xmlns:w=\"http://schemas.openxmlformats.org/wordprocessingml/2006/main\"><w:body><w:p><w:pPr><w:rPr><w:b/></w:rPr></w:pPr><w:r><w:rPr><w:b/></w:rPr><w:t>Git Dersi \u00c7al\u0131\u015fmalar\u0131</w:t></w:r></w:p></w:body></w:document></pkg:xmlData></pkg:part></pkg:package>'\n$titleRange.InsertXML($titleOoxml)\n\n# --- 2) Append two new paragraphs after the existing \"Belgemin ilk\n#        paragraf\u0131\" paragraph: \"\u0130kinci paragraf\u0131\" and \"\u00dc\u00e7\u00fcnc\u00fc paragraf\u0131\". ---\n$lastRange = $d.Paragraphs($d.Paragraphs.Count).Range\n$lastRange.InsertParagraphAfter()\n$lastRange.Collapse(0)\n$d.Paragraphs($d.Paragraphs.Count).Range.Text = \"\u0130kinci paragraf\u0131\"\n\n$lastRange2 = $d.Paragraphs($d.Paragraphs.Count).Range\n$lastRange2.InsertParagraphAfter()\n$lastRange2.Collapse(0)\n$d.Paragraphs($d.Paragraphs.Count).Range.Text = \"\u00dc\u00e7\u00fcnc\u00fc paragraf\u0131\"\n"}
